$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(667.4, 1332.6, 271, 1729, 282.9, 1717.1, 320.9, 1679.1, 1931.5, 68.5, 1718.8, 281.2, 1140, 860)
    3  = @(1288.8, 711.2, 1329.5, 670.5, 1859.9, 140.1, 1175.4, 824.6, 438.4, 1561.6, 285.7, 1714.3, 479.3, 1520.7)
    4  = @(1950.5, 49.5, 681.3, 1318.7, 596.1, 1403.9, 929.8, 1070.2, 649.8, 1350.2, 527.9, 1472.1, 1603.4, 396.6)
    5  = @(64.40000000000001, 1935.6, 1747.6, 252.4, 1050.1, 949.9, 1565.6, 434.4, 1219.1, 780.9, 1461, 539, 762.1, 1237.9)
    6  = @(382.7, 1617.3, 858.3, 1141.7, 1525.3, 474.7, 1857, 143, 1359.5, 640.5, 1874.7, 125.3, 1863, 137)
    7  = @(1511, 489, 1792, 208, 71.2, 1928.8, 647.2, 1352.8, 759, 1241, 182.8, 1817.2, 521.6, 1478.4)
    8  = @(1220.5, 779.5, 189.5, 1810.5, 1338.1, 661.9, 1403.4, 596.6, 78.8, 1921.2, 878.1, 1121.9, 1395.6, 604.4)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Count; $i++) {
        # Column B is index 2 in 1-based column numbering
        $col = $i + 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
